# This script applies a row-wise shuffle of the data rows (2-43) in the
# "Haba" sheet. The columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), O (Origen) and
# P (Precio $/Kg) are permuted across rows according to the mapping
# below (target row -> source row, both referring to the *original*
# values before this script runs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (source row refers to the original/before values)
$rowMap = @{
  2  = 42
  3  = 22
  4  = 24
  5  = 33
  6  = 10
  7  = 15
  8  = 4
  9  = 13
  10 = 27
  11 = 31
  12 = 14
  13 = 6
  14 = 19
  15 = 29
  16 = 26
  17 = 17
  18 = 18
  19 = 20
  20 = 28
  21 = 39
  22 = 2
  23 = 12
  24 = 38
  25 = 3
  26 = 35
  27 = 40
  28 = 23
  29 = 5
  30 = 37
  31 = 9
  32 = 7
  33 = 8
  34 = 21
  35 = 34
  36 = 30
  37 = 43
  38 = 32
  39 = 41
  40 = 16
  41 = 25
  42 = 36
  43 = 11
}

# Columns that move together as a group for each data row.
$cols = @("D", "J", "K", "L", "M", "O", "P")

# Snapshot the original values of every affected column for rows 2..43
# before any writes happen, so the permutation is applied consistently.
$snapshot = @{}
for ($r = 2; $r -le 43; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Write the new values for every row using the snapshot of the source row.
for ($r = 2; $r -le 43; $r++) {
    $srcRow = $rowMap[$r]
    if ($srcRow -eq $r) {
        continue
    }
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $srcVals[$c]
    }
}
